# Auto-generated script to update resultados_accd_01.xlsx
# Updates columns B (recommended_rank tuple, shared across all sheets) and C (ideal_rank tuple, per-sheet) on all 7 sheets
$wb = $excel.ActiveWorkbook

# Phase 1: recommended_rank column (B2:B8) - identical values on every sheet
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 2).Value = '(''DecisionTree'', (2.233954249025107, 0.7540653475337438))'
$ws.Cells.Item(3, 2).Value = '(''Knn10'', (1.479402190895932, 0.7319053522713357))'
$ws.Cells.Item(4, 2).Value = '(''LDA'', (1.4081045806319594, 0.6753846338630222))'
$ws.Cells.Item(5, 2).Value = '(''Knn5'', (1.3108624895035417, 0.7128458097177036))'
$ws.Cells.Item(6, 2).Value = '(''Knn1'', (0.8463931672647643, 0.6472733703791613))'
$ws.Cells.Item(7, 2).Value = '(''RandomForest'', (0.8375589217023194, 0.7594959733070695))'
$ws.Cells.Item(8, 2).Value = '(''NaiveBayes'', (0.3523928679837268, 0.5294468484269127))'

# Phase 2: ideal_rank column (C2:C8) - per-sheet values
# Sheet 1: abalone
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 3).Value = '(''DecisionTree'', (1.317034368673276, 0.7540653475337438))'
$ws.Cells.Item(3, 3).Value = '(''Knn10'', (1.0335739632270662, 0.7319053522713357))'
$ws.Cells.Item(4, 3).Value = '(''Knn5'', (1.0052177117125636, 0.7128458097177036))'
$ws.Cells.Item(5, 3).Value = '(''LDA'', (0.9552120772111133, 0.6753846338630222))'
$ws.Cells.Item(6, 3).Value = '(''RandomForest'', (0.9169749356503704, 0.7594959733070695))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.8786059638221791, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''NaiveBayes'', (0.8481621349153448, 0.5294468484269127))'

# Sheet 2: adult
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 3).Value = '(''DecisionTree'', (1.2928406096289797, 0.7540653475337438))'
$ws.Cells.Item(3, 3).Value = '(''Knn10'', (1.101322324033936, 0.7319053522713357))'
$ws.Cells.Item(4, 3).Value = '(''Knn5'', (1.0612108946257937, 0.7128458097177036))'
$ws.Cells.Item(5, 3).Value = '(''LDA'', (1.0276330403700729, 0.6753846338630222))'
$ws.Cells.Item(6, 3).Value = '(''RandomForest'', (0.926484135829032, 0.7594959733070695))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.910468185520431, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''NaiveBayes'', (0.7224559695684375, 0.5294468484269127))'

# Sheet 3: banknote
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 3).Value = '(''DecisionTree'', (1.3040207515107816, 0.7540653475337438))'
$ws.Cells.Item(3, 3).Value = '(''Knn10'', (1.0767032399724972, 0.7319053522713357))'
$ws.Cells.Item(4, 3).Value = '(''Knn5'', (1.0364286097620394, 0.7128458097177036))'
$ws.Cells.Item(5, 3).Value = '(''LDA'', (1.0289667489904193, 0.6753846338630222))'
$ws.Cells.Item(6, 3).Value = '(''RandomForest'', (0.9440752102657682, 0.7594959733070695))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.8781849628288728, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''NaiveBayes'', (0.7873294306442793, 0.5294468484269127))'

# Sheet 4: car
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 3).Value = '(''DecisionTree'', (1.2818136118202792, 0.7540653475337438))'
$ws.Cells.Item(3, 3).Value = '(''Knn10'', (1.0864831418368561, 0.7319053522713357))'
$ws.Cells.Item(4, 3).Value = '(''Knn5'', (1.0543798493069336, 0.7128458097177036))'
$ws.Cells.Item(5, 3).Value = '(''LDA'', (1.0362931273300686, 0.6753846338630222))'
$ws.Cells.Item(6, 3).Value = '(''RandomForest'', (0.9289772820964696, 0.7594959733070695))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.9140891809248314, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''NaiveBayes'', (0.7327397044023539, 0.5294468484269127))'

# Sheet 5: chess1
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 3).Value = '(''DecisionTree'', (1.2643666747871574, 0.7540653475337438))'
$ws.Cells.Item(3, 3).Value = '(''Knn10'', (1.092438029443945, 0.7319053522713357))'
$ws.Cells.Item(4, 3).Value = '(''Knn5'', (1.0477354130651249, 0.7128458097177036))'
$ws.Cells.Item(5, 3).Value = '(''LDA'', (1.0323956944926018, 0.6753846338630222))'
$ws.Cells.Item(6, 3).Value = '(''RandomForest'', (0.9144641763688443, 0.7594959733070695))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.899253550702905, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''NaiveBayes'', (0.8039147510219499, 0.5294468484269127))'

# Sheet 6: chess2
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 3).Value = '(''LDA'', (1.0424900022499681, 0.6753846338630222))'
$ws.Cells.Item(3, 3).Value = '(''DecisionTree'', (1.001835319538069, 0.7540653475337438))'
$ws.Cells.Item(4, 3).Value = '(''Knn10'', (0.9380387414990355, 0.7319053522713357))'
$ws.Cells.Item(5, 3).Value = '(''Knn5'', (0.9013388739689342, 0.7128458097177036))'
$ws.Cells.Item(6, 3).Value = '(''NaiveBayes'', (0.8991628017162682, 0.5294468484269127))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.8118854646461315, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''RandomForest'', (0.764657542026064, 0.7594959733070695))'

# Sheet 7: contraceptive
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 3).Value = '(''DecisionTree'', (1.3130942999717679, 0.7540653475337438))'
$ws.Cells.Item(3, 3).Value = '(''Knn10'', (1.0722123898583404, 0.7319053522713357))'
$ws.Cells.Item(4, 3).Value = '(''Knn5'', (1.031972445253473, 0.7128458097177036))'
$ws.Cells.Item(5, 3).Value = '(''LDA'', (1.031195942505484, 0.6753846338630222))'
$ws.Cells.Item(6, 3).Value = '(''RandomForest'', (0.93673667049728, 0.7594959733070695))'
$ws.Cells.Item(7, 3).Value = '(''Knn1'', (0.896225895899366, 0.6472733703791613))'
$ws.Cells.Item(8, 3).Value = '(''NaiveBayes'', (0.7701337078425116, 0.5294468484269127))'

